$wb = $excel.ActiveWorkbook

# --- optimization_parameters sheet: insert an "L_curve" parameter row and
# rename the "Model" label to "production_function" ------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Insert a new row above the old row 9 ("estimate_params") so everything
# below shifts down by one, exactly like Excel's Insert Row command.
$ws.Rows(9).Insert()

# The old row 8 ("Model" / "MM") stays in place but its label changes.
$ws.Range("A8").Value = "production_function"

# The freshly inserted row 9 becomes the new "L_curve" parameter (value 0).
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

# Make this the active sheet/tab and select the newly inserted row, matching
# how a user would have just finished typing the new parameter in.
$ws.Activate()
$ws.Range("A9:B9").Select()
